$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 392, shifting existing rows 392:450 down to 393:451
$ws.Rows.Item(392).Insert()

# Populate the newly inserted row 392 with the new data record
$ws.Range("A392").Value = 1
$ws.Range("B392").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C392").Value = "Arica y Parinacota"
$ws.Range("D392").Value = 44984
$ws.Range("E392").Value = 15
$ws.Range("F392").Value = 100114013
$ws.Range("G392").Value = "Zanahoria"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 70
$ws.Range("K392").Value = 9000
$ws.Range("L392").Value = 10000
$ws.Range("M392").Value = 9500
$ws.Range("N392").Value = "$/saco 25 kilos"
$ws.Range("O392").Value = "Región de Arica y Parinacota"
$ws.Range("P392").Value = 380
$ws.Range("Q392").Value = 25
$ws.Range("R392").Value = "Hortaliza"
